$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11 (anchor G11=5533)
$ws.Range("H11").Value = 4620
$ws.Range("I11").Value = 4620
$ws.Range("K11").Value = 4620
$ws.Range("M11").Value = -4480
# Row 51 (anchor G51=5486)
$ws.Range("H51").Value = 6649.1665
$ws.Range("J51").Value = 6998.3335
$ws.Range("L51").Value = 6998.3335
$ws.Range("N51").Value = -7966.3335
# Row 80 (anchor G80=12605)
$ws.Range("H80").Value = 1373.8948
$ws.Range("I80").Value = 1548.125
$ws.Range("J80").Value = 1247.1818
$ws.Range("K80").Value = 4644.375
$ws.Range("L80").Value = 3741.5454
$ws.Range("M80").Value = -3646.375
$ws.Range("N80").Value = -5737.5454
# Row 83 (anchor G83=12605)
$ws.Range("H83").Value = 1373.8948
$ws.Range("I83").Value = 1548.125
$ws.Range("J83").Value = 1247.1818
$ws.Range("K83").Value = 13933.125
$ws.Range("L83").Value = 11224.6362
$ws.Range("M83").Value = -8941.125
$ws.Range("N83").Value = -21208.6362
# Row 88 (anchor G88=12608)
$ws.Range("H88").Value = 2880.5
$ws.Range("I88").Value = 2328.3333
$ws.Range("J88").Value = 3432.6667
$ws.Range("K88").Value = 2328.3333
$ws.Range("L88").Value = 3432.6667
$ws.Range("M88").Value = -1922.3333
$ws.Range("N88").Value = -4244.6667
# Row 91 (anchor G91=12608)
$ws.Range("H91").Value = 2880.5
$ws.Range("I91").Value = 2328.3333
$ws.Range("J91").Value = 3432.6667
$ws.Range("K91").Value = 2328.3333
$ws.Range("L91").Value = 3432.6667
$ws.Range("M91").Value = -924.3332999999998
$ws.Range("N91").Value = -6240.6667
# Row 111 (anchor G111=27768)
$ws.Range("H111").Value = 2750
$ws.Range("I111").Value = 3000
$ws.Range("J111").Value = 2500
$ws.Range("K111").Value = 9000
$ws.Range("L111").Value = 7500
$ws.Range("M111").Value = -5933
$ws.Range("N111").Value = -13634
# Row 137 (anchor G137=44013)
$ws.Range("H137").Value = 2068.2856
$ws.Range("J137").Value = 1747.5
$ws.Range("L137").Value = 5242.5
$ws.Range("N137").Value = -10342.5
# Row 138 (anchor G138=44169)
$ws.Range("H138").Value = 3297.75
$ws.Range("I138").Value = 1998.5
$ws.Range("J138").Value = 4597
$ws.Range("K138").Value = 5995.5
$ws.Range("L138").Value = 13791
$ws.Range("M138").Value = -855.5
$ws.Range("N138").Value = -24071

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5 (anchor G5=5091)
$ws.Range("H5").Value = 1174.8182
$ws.Range("I5").Value = 2221.8
$ws.Range("K5").Value = 2221.8
$ws.Range("M5").Value = -2109.8
# Row 61 (anchor G61=43999)
$ws.Range("H61").Value = 1500
$ws.Range("I61").Value = 1500
$ws.Range("K61").Value = 1500
$ws.Range("M61").Value = -1288
# Row 88 (anchor G88=12530)
$ws.Range("H88").Value = 3125
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
# Row 91 (anchor G91=12530)
$ws.Range("H91").Value = 3125
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
# Row 110 (anchor G110=27708)
$ws.Range("H110").Value = 6204.7144
$ws.Range("I110").Value = 6204.7144
$ws.Range("K110").Value = 6204.7144
$ws.Range("M110").Value = -4159.7144
# Row 132 (anchor G132=43997)
$ws.Range("H132").Value = 3454.647
$ws.Range("I132").Value = 1973
$ws.Range("K132").Value = 5919
$ws.Range("M132").Value = -3389
# Row 136 (anchor G136=43999)
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("M136").Value = -1950

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4 (anchor G4=5091)
$ws.Range("H4").Value = 1174.8182
$ws.Range("I4").Value = 2221.8
$ws.Range("K4").Value = 2221.8
$ws.Range("M4").Value = -2106.8
# Row 86 (anchor G86=12526)
$ws.Range("H86").Value = 1908
$ws.Range("I86").Value = 2012
$ws.Range("J86").Value = 1492
$ws.Range("K86").Value = 2012
$ws.Range("L86").Value = 1492
$ws.Range("M86").Value = -889
$ws.Range("N86").Value = -3738
# Row 89 (anchor G89=12526)
$ws.Range("H89").Value = 1908
$ws.Range("I89").Value = 2012
$ws.Range("J89").Value = 1492
$ws.Range("K89").Value = 10060
$ws.Range("L89").Value = 7460
$ws.Range("M89").Value = -4444
$ws.Range("N89").Value = -18692

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 17 (anchor G17=1823)
$ws.Range("H17").Value = 9950
$ws.Range("J17").Value = 9950
$ws.Range("L17").Value = 9950
$ws.Range("N17").Value = -10298
# Row 25 (anchor G25=1895)
$ws.Range("H25").Value = 2368
$ws.Range("I25").Value = 2310
$ws.Range("J25").Value = 2513
$ws.Range("K25").Value = 2310
$ws.Range("L25").Value = 2513
$ws.Range("M25").Value = -2136
$ws.Range("N25").Value = -2861
# Row 31 (anchor G31=44023)
$ws.Range("H31").Value = 1704.1666
$ws.Range("I31").Value = 1895
$ws.Range("J31").Value = 750
$ws.Range("K31").Value = 1895
$ws.Range("L31").Value = 750
$ws.Range("M31").Value = -1600
$ws.Range("N31").Value = -1340
# Row 34 (anchor G34=44023)
$ws.Range("H34").Value = 1704.1666
$ws.Range("I34").Value = 1895
$ws.Range("J34").Value = 750
$ws.Range("K34").Value = 1895
$ws.Range("L34").Value = 750
$ws.Range("M34").Value = -1693
$ws.Range("N34").Value = -1154
# Row 41 (anchor G41=1917)
$ws.Range("H41").Value = 8341.6
$ws.Range("I41").Value = 8341.6
$ws.Range("K41").Value = 8341.6
$ws.Range("M41").Value = -7913.6
# Row 50 (anchor G50=1862)
$ws.Range("H50").Value = 9092
$ws.Range("J50").Value = 9092
$ws.Range("L50").Value = 9092
$ws.Range("N50").Value = -10342
# Row 51 (anchor G51=2039)
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
# Row 54 (anchor G54=2413)
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
# Row 55 (anchor G55=1855)
$ws.Range("H55").Value = 10073
$ws.Range("I55").Value = 10073
$ws.Range("K55").Value = 10073
$ws.Range("M55").Value = -9758
# Row 61 (anchor G61=2039)
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
# Row 69 (anchor G69=11911)
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
# Row 72 (anchor G72=11911)
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
# Row 132 (anchor G132=44019)
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12 (anchor G12=4854)
$ws.Range("H12").Value = 75.333336
$ws.Range("I12").Value = 75.75
$ws.Range("K12").Value = 227.25
$ws.Range("M12").Value = -54.25
# Row 132 (anchor G132=43972)
$ws.Range("H132").Value = 1996
$ws.Range("J132").Value = 1996
$ws.Range("L132").Value = 17964
$ws.Range("N132").Value = -23024
# Row 139 (anchor G139=44102)
$ws.Range("H139").Value = 3166
$ws.Range("I139").Value = 3298.8
$ws.Range("K139").Value = 9896.400000000001
$ws.Range("M139").Value = -4756.400000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122 (anchor G122=36182)
$ws.Range("H122").Value = 2176.5557
$ws.Range("I122").Value = 1689.1364
$ws.Range("J122").Value = 4321.2
$ws.Range("K122").Value = 5067.4092
$ws.Range("L122").Value = 12963.6
$ws.Range("M122").Value = -2617.4092
$ws.Range("N122").Value = -17863.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40 (anchor G40=36248)
$ws.Range("H40").Value = 3061.7896
$ws.Range("J40").Value = 1517.25
$ws.Range("L40").Value = 1517.25
$ws.Range("N40").Value = -1789.25
# Row 61 (anchor G61=27740)
$ws.Range("H61").Value = 1780.4736
$ws.Range("I61").Value = 1616.0714
$ws.Range("K61").Value = 1616.0714
$ws.Range("M61").Value = -1414.0714
# Row 93 (anchor G93=19993)
$ws.Range("H93").Value = 2260
$ws.Range("I93").Value = 1450
$ws.Range("J93").Value = 5500
$ws.Range("K93").Value = 1450
$ws.Range("L93").Value = 5500
$ws.Range("M93").Value = -202
$ws.Range("N93").Value = -7996
# Row 113 (anchor G113=27740)
$ws.Range("H113").Value = 1780.4736
$ws.Range("I113").Value = 1616.0714
$ws.Range("K113").Value = 1616.0714
$ws.Range("M113").Value = 553.9286
